$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Add new row 18 (new IPA test case), mirroring the styling pattern of row 17
$ws.Range("A18").Value = "IPAIAM0057"
$ws.Range("B18").Value = "OPQA-4525||OPQA-4526||OPQA-4527"
$ws.Range("C18").Value = "Verify that the STeAM Step Up Auth Modal should be presented to the user without a pre-populated email address when user has a valid Neon session token and is navigating within the same browser window."
$ws.Range("D18").Value = "Y"

$ws.Range("A18").Style = $ws.Range("A17").Style
$ws.Range("B18").Style = $ws.Range("B17").Style
$ws.Range("C18").Style = $ws.Range("A17").Style
$ws.Range("D18").Style = $ws.Range("D17").Style
$ws.Range("E18").Style = $ws.Range("E17").Style

$ws.Rows.Item(18).RowHeight = 45

# Update the current selection/scroll position to reflect the new last row
$ws.Application.ActiveWindow.ScrollRow = 16
$ws.Range("D19").Select()
